$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row restructuring -------------------------------------------------
# The "nakit_orani" row (row 3) was removed entirely, shifting the rows
# below it up by one; a new row ("max indirimli kredi miktari") was then
# appended at the bottom.
$ws.Rows(3).Delete()

# --- Updated values ------------------------------------------------------
$ws.Range("B2").Value = 0.5        # kredi_orani
$ws.Range("B5").Value = 0.09       # indirimli_yillik_faiz

# --- New row at the bottom ------------------------------------------------
$ws.Range("A6").Value = "max indirimli kredi miktari"
$ws.Range("B6").Value = 100000

# --- Formatting ------------------------------------------------------------
# Bold header row
$ws.Range("A1:B1").Font.Bold = $true

# Percent style for the interest-rate rows
$ws.Range("B4:B5").Style = "Percent"

# Comma style (thousands separator, no decimals) for the new max-credit row
$ws.Range("B6").Style = "Comma"
$ws.Range("B6").NumberFormat = '_-* #,##0_-;\-* #,##0_-;_-* "-"??_-;_-@_-'

# --- Column widths ---------------------------------------------------------
$ws.Columns("A").ColumnWidth = 24
$ws.Columns("B").ColumnWidth = 10.7

# --- Selection ---------------------------------------------------------------
$ws.Range("D11").Select()

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
